$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 2.77
$ws.Range("J2").Value = 2.82
$ws.Range("K2").Value = 2.22
$ws.Range("L2").Value = 3.25
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 8
$ws.Range("P2").Value = 3.65
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2.02
$ws.Range("S2").Value = 1.34
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 1.6
$ws.Range("W2").Value = 9.25
$ws.Range("X2").Value = 12
$ws.Range("Y2").Value = 9
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 17.5
$ws.Range("AB2").Value = 24
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 6.8
$ws.Range("AH2").Value = 10.75
$ws.Range("AI2").Value = 16
$ws.Range("AJ2").Value = 10.25
$ws.Range("AK2").Value = 35
$ws.Range("AL2").Value = 22
$ws.Range("AM2").Value = 26
$ws.Range("AN2").Value = 4.4
$ws.Range("AO2").Value = 11.75
$ws.Range("AP2").Value = 18
$ws.Range("AQ2").Value = 45
$ws.Range("AR2").Value = 70
$ws.Range("AT2").Value = 3
$ws.Range("AU2").Value = 6.6
$ws.Range("AW2").Value = 4.9
$ws.Range("AX2").Value = 14.5
$ws.Range("AY2").Value = 19.5
$ws.Range("AZ2").Value = 60
$ws.Range("BA2").Value = 80
